$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D(4), M(13), N(14), O(15), P(16), S(19).
# Row 2-41 data: Fecha(D), Volumen(M), Precio minimo(N), Precio maximo(O), Precio promedio ponderado(P), Precio $/Kg(S)
$data = @(
    @(44431, 160, 21000, 22000, 21500, 1075),
    @(44467, 200, 20000, 21000, 20500, 1025),
    @(44879, 100, 28000, 30000, 29000, 1450),
    @(44448, 100, 20000, 21000, 20500, 1025),
    @(44315, 100, 20000, 21000, 20500, 1025),
    @(44365, 100, 20000, 21000, 20500, 1025),
    @(44782, 200, 23500, 24000, 23750, 1188),
    @(44427, 200, 20000, 21000, 20500, 1025),
    @(44435, 260, 20000, 22000, 21115, 1056),
    @(44333, 100, 19500, 20000, 19750, 988),
    @(44874, 240, 29000, 30000, 29500, 1475),
    @(44473, 40, 19500, 20000, 19750, 988),
    @(44326, 160, 19500, 20000, 19750, 988),
    @(44335, 200, 19000, 20000, 19500, 975),
    @(44882, 120, 28000, 30000, 29000, 1450),
    @(44301, 100, 18000, 19000, 18500, 925),
    @(44462, 100, 19500, 20000, 19750, 988),
    @(44880, 100, 28000, 30000, 29000, 1450),
    @(44442, 140, 20000, 21000, 20500, 1025),
    @(44420, 160, 20000, 21000, 20500, 1025),
    @(44410, 200, 20000, 21000, 20500, 1025),
    @(44343, 100, 19500, 20000, 19750, 988),
    @(44364, 140, 20000, 21000, 20500, 1025),
    @(44810, 100, 27000, 28000, 27500, 1375),
    @(44428, 100, 20000, 21000, 20500, 1025),
    @(44784, 160, 27000, 28000, 27500, 1375),
    @(44445, 160, 20000, 21000, 20500, 1025),
    @(44778, 100, 23000, 24000, 23500, 1175),
    @(44350, 160, 19000, 20000, 19500, 975),
    @(44809, 60, 27000, 28000, 27500, 1375),
    @(44781, 160, 23000, 24000, 23500, 1175),
    @(44434, 100, 20000, 21000, 20500, 1025),
    @(44407, 160, 20000, 21000, 20500, 1025),
    @(44418, 200, 20000, 21000, 20500, 1025),
    @(44474, 200, 19000, 20000, 19500, 975),
    @(44336, 100, 19500, 20000, 19750, 988),
    @(44776, 160, 23000, 24000, 23500, 1175),
    @(44417, 160, 20000, 21000, 20500, 1025),
    @(44466, 100, 20000, 21000, 20500, 1025),
    @(44441, 160, 20000, 21000, 20500, 1025)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals[5]   # S - Precio $/Kg
}

Write-Host "Done updating rows 2-41"